$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

# Force the whole new row to a text number format first so that
# numeric-looking strings (phone number) are kept as text and so that
# empty-string assignments still leave a (empty) cell behind, matching
# the sibling rows in the log which always have one <c> per column.
$rowRange = $ws.Range("A$row`:G$row")
$rowRange.NumberFormat = "@"

$ws.Range("A$row").Value = "2025-06-17T02:03:43.065517"
$ws.Range("B$row").Value = "Akash"
$ws.Range("C$row").Value = ""
$ws.Range("D$row").Value = "9386776355"
$ws.Range("E$row").Value = ""
$ws.Range("F$row").Value = "Unable to get second tranche"
$ws.Range("G$row").Value = "Sonpur"

# The other data rows use the default (unstyled) cell format, so put the
# style back to Normal now that the values are safely typed as text.
$rowRange.Style = "Normal"
